$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

$sub5 = [char]0x2085

Set-TextValue $ws.Range("D2") "25.840.86"
Set-TextValue $ws.Range("E2") "  -0.25%  "
Set-TextValue $ws.Range("D3") "1.637.48"
Set-TextValue $ws.Range("E3") "  -0.02%  "
Set-TextValue $ws.Range("E4") "  +0.01%  "
Set-TextValue $ws.Range("D5") "216.05"
Set-TextValue $ws.Range("E5") "  +0.67%  "
Set-TextValue $ws.Range("D6") "0.5077"
Set-TextValue $ws.Range("E6") "  +0.17%  "
Set-TextValue $ws.Range("E7") "  +0.13%  "
Set-TextValue $ws.Range("D8") "0.2583"
Set-TextValue $ws.Range("E8") "  +0.34%  "
Set-TextValue $ws.Range("D9") "0.06443"
Set-TextValue $ws.Range("E9") "  +1.32%  "
Set-TextValue $ws.Range("D10") "19.59"
Set-TextValue $ws.Range("E10") "  -1.03%  "
Set-TextValue $ws.Range("D11") "0.07789"
Set-TextValue $ws.Range("E11") "  +0.58%  "
Set-TextValue $ws.Range("D12") "4.286"
Set-TextValue $ws.Range("E12") "  -0.24%  "
Set-TextValue $ws.Range("D13") "1.863.55"
Set-TextValue $ws.Range("E13") "  -0.02%  "
Set-TextValue $ws.Range("D14") "1.633.56"
Set-TextValue $ws.Range("E14") "  -0.24%  "
Set-TextValue $ws.Range("D15") "0.5635"
Set-TextValue $ws.Range("E15") "  +3.11%  "
Set-TextValue $ws.Range("D16") "0.0${sub5}7616"
Set-TextValue $ws.Range("E16") "  -1.68%  "
Set-TextValue $ws.Range("E17") "  -1.39%  "
Set-TextValue $ws.Range("D18") "25.859.46"
Set-TextValue $ws.Range("E19") "  +0.12%  "
Set-TextValue $ws.Range("D20") "195.68"
Set-TextValue $ws.Range("E20") "  -0.05%  "
Set-TextValue $ws.Range("D21") "4.329"
Set-TextValue $ws.Range("E21") "  -2.90%  "
Set-TextValue $ws.Range("D22") "9.892"
Set-TextValue $ws.Range("E22") "  -0.52%  "
Set-TextValue $ws.Range("D23") "6.098"
Set-TextValue $ws.Range("E23") "  -0.62%  "
Set-TextValue $ws.Range("E24") "  +0.02%  "
Set-TextValue $ws.Range("D25") "1.795"
Set-TextValue $ws.Range("E25") "  -5.52%  "
Set-TextValue $ws.Range("D26") "0.1270"
Set-TextValue $ws.Range("E26") "  +0.93%  "
Set-TextValue $ws.Range("D27") "139.92"
Set-TextValue $ws.Range("E27") "  -2.37%  "
Set-TextValue $ws.Range("E28") "  -0.67%  "
Set-TextValue $ws.Range("D29") "15.47"
Set-TextValue $ws.Range("E29") "  -1.15%  "
Set-TextValue $ws.Range("E30") "  +0.49%  "
Set-TextValue $ws.Range("D31") "0.04885"
Set-TextValue $ws.Range("E31") "  +0.15%  "
Set-TextValue $ws.Range("D32") "3.299"
Set-TextValue $ws.Range("E32") "  +1.76%  "
Set-TextValue $ws.Range("D33") "3.231"
Set-TextValue $ws.Range("E33") "  +0.98%  "
Set-TextValue $ws.Range("E34") "  +0.18%  "
Set-TextValue $ws.Range("D35") "2.368"
Set-TextValue $ws.Range("E35") "  -0.38%  "
Set-TextValue $ws.Range("D36") "0.9036"
Set-TextValue $ws.Range("E36") "  -1.19%  "
Set-TextValue $ws.Range("D37") "2.575"
Set-TextValue $ws.Range("E37") "  +0.19%  "
Set-TextValue $ws.Range("D38") "1.130.47"
Set-TextValue $ws.Range("E38") "  -0.05%  "
Set-TextValue $ws.Range("D39") "0.5505"
Set-TextValue $ws.Range("E39") "  -0.30%  "
Set-TextValue $ws.Range("D40") "0.01564"
Set-TextValue $ws.Range("E40") "  -0.06%  "
Set-TextValue $ws.Range("E41") "  -0.69%  "
Set-TextValue $ws.Range("D42") "5.532"
Set-TextValue $ws.Range("E42") "  -1.18%  "
Set-TextValue $ws.Range("D43") "0.8009"
Set-TextValue $ws.Range("E43") "  -0.44%  "
Set-TextValue $ws.Range("E44") "  -0.77%  "
Set-TextValue $ws.Range("D45") "1.773.86"
Set-TextValue $ws.Range("E45") "  +0.02%  "
Set-TextValue $ws.Range("E46") "  -6.80%  "
Set-TextValue $ws.Range("E47") "  -1.01%  "
Set-TextValue $ws.Range("D48") "55.44"
Set-TextValue $ws.Range("E48") "  +0.38%  "
Set-TextValue $ws.Range("D49") "7.705"
Set-TextValue $ws.Range("D50") "0.05053"
Set-TextValue $ws.Range("E50") "  -2.48%  "
Set-TextValue $ws.Range("D51") "1.003"
Set-TextValue $ws.Range("E51") "  +0.09%  "
